$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the three target duration values from 0.1 to 0.01
$ws.Range("C14").Value = 0.01
$ws.Range("C20").Value = 0.01
$ws.Range("C42").Value = 0.01

# Update the view to match new scroll position / zoom level
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.Zoom = 161
